$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: RandomForestRegressor - values change, name unchanged
$ws.Range("B3").Value = 0.978416500096945
$ws.Range("C3").Value = 0.9767452644274428
$ws.Range("D3").Value = 0.9767038059856029

# Row 4: GradientBoostingRegressor -> DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9793034381245443
$ws.Range("C4").Value = 0.9802699348530096
$ws.Range("D4").Value = 0.980289892381167

# Row 5: AdaBoostRegressor -> MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.9763731090403213
$ws.Range("C5").Value = 0.9805226687536566
$ws.Range("D5").Value = 0.981632033308304
